$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1516.3125
$ws.Range("J17").Value = 1564.0667
$ws.Range("L17").Value = 4692.2001
$ws.Range("N17").Value = -5028.2001

$ws.Range("H33").Value = 80359.2
$ws.Range("I33").Value = 125236.875
$ws.Range("J33").Value = 576.6667
$ws.Range("K33").Value = 125236.875
$ws.Range("L33").Value = 576.6667
$ws.Range("M33").Value = -125007.875
$ws.Range("N33").Value = -1034.6667

$ws.Range("H41").Value = 2146.2
$ws.Range("I41").Value = 2677.3333
$ws.Range("J41").Value = 1349.5
$ws.Range("K41").Value = 2677.3333
$ws.Range("L41").Value = 1349.5
$ws.Range("M41").Value = -2237.3333
$ws.Range("N41").Value = -2229.5

$ws.Range("H92").Value = 445.8
$ws.Range("I92").Value = 447.75
$ws.Range("J92").Value = 442.875
$ws.Range("K92").Value = 447.75
$ws.Range("L92").Value = 442.875
$ws.Range("M92").Value = 800.25
$ws.Range("N92").Value = -2938.875

$ws.Range("H99").Value = 1121.1111
$ws.Range("I99").Value = 258.2
$ws.Range("K99").Value = 774.5999999999999
$ws.Range("M99").Value = 723.4000000000001

$ws.Range("H115").Value = 4777.6
$ws.Range("I115").Value = 4777.6
$ws.Range("K115").Value = 14332.8
$ws.Range("M115").Value = -12765.8

$ws.Range("H125").Value = 994.3333
$ws.Range("J125").Value = 1071.7778
$ws.Range("L125").Value = 9646.0002
$ws.Range("N125").Value = -14566.0002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 34682.668
$ws.Range("J44").Value = 34682.668
$ws.Range("L44").Value = 34682.668
$ws.Range("N44").Value = -35658.668

$ws.Range("H55").Value = 30693.166
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 30693.166
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 30693.166
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -31323.166

$ws.Range("H74").Value = 1935.2812
$ws.Range("I74").Value = 2062.7273
$ws.Range("K74").Value = 2062.7273
$ws.Range("M74").Value = -1188.7273

$ws.Range("H77").Value = 1935.2812
$ws.Range("I77").Value = 2062.7273
$ws.Range("K77").Value = 10313.6365
$ws.Range("M77").Value = -5945.636500000001

$ws.Range("H97").Value = 2007.8096
$ws.Range("I97").Value = 1529.125
$ws.Range("J97").Value = 3539.6
$ws.Range("K97").Value = 1529.125
$ws.Range("L97").Value = 3539.6
$ws.Range("M97").Value = -1033.125
$ws.Range("N97").Value = -4531.6

$ws.Range("H102").Value = 1543.2307
$ws.Range("I102").Value = 1527.75
$ws.Range("J102").Value = 1568
$ws.Range("K102").Value = 1527.75
$ws.Range("L102").Value = 1568
$ws.Range("M102").Value = 94.25
$ws.Range("N102").Value = -4812

$ws.Range("H110").Value = 2442.2222
$ws.Range("I110").Value = 1501
$ws.Range("K110").Value = 1501
$ws.Range("M110").Value = 544

$ws.Range("H131").Value = 35749.25
$ws.Range("I131").Value = 18000
$ws.Range("J131").Value = 41665.668
$ws.Range("K131").Value = 18000
$ws.Range("L131").Value = 41665.668
$ws.Range("N131").Value = -51745.668
$ws.Range("M131").Value = -12960

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 1534
$ws.Range("I64").Value = 1002
$ws.Range("J64").Value = 1800
$ws.Range("K64").Value = 1002
$ws.Range("L64").Value = 1800
$ws.Range("M64").Value = -777
$ws.Range("N64").Value = -2250

$ws.Range("H67").Value = 1534
$ws.Range("I67").Value = 1002
$ws.Range("J67").Value = 1800
$ws.Range("K67").Value = 1002
$ws.Range("L67").Value = 1800
$ws.Range("M67").Value = -222
$ws.Range("N67").Value = -3360

$ws.Range("H94").Value = 1560.3334
$ws.Range("I94").Value = 1427.6
$ws.Range("J94").Value = 1726.25
$ws.Range("K94").Value = 1427.6
$ws.Range("L94").Value = 1726.25
$ws.Range("M94").Value = -976.5999999999999
$ws.Range("N94").Value = -2628.25

$ws.Range("H99").Value = 1564.5
$ws.Range("I99").Value = 1869.6666
$ws.Range("K99").Value = 1869.6666
$ws.Range("M99").Value = -371.6666

$ws.Range("H140").Value = 147497.5
$ws.Range("J140").Value = 147497.5
$ws.Range("L140").Value = 147497.5
$ws.Range("N140").Value = -157857.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1703.3334
$ws.Range("J94").Value = 1818.6666
$ws.Range("L94").Value = 1818.6666
$ws.Range("N94").Value = -2720.6666

$ws.Range("H105").Value = 3096.5715
$ws.Range("I105").Value = 1654.2858
$ws.Range("J105").Value = 4538.857
$ws.Range("K105").Value = 1654.2858
$ws.Range("L105").Value = 4538.857
$ws.Range("M105").Value = 92.71419999999989
$ws.Range("N105").Value = -8032.857

$ws.Range("H134").Value = 1807.0344
$ws.Range("I134").Value = 1631
$ws.Range("K134").Value = 4893
$ws.Range("M134").Value = -2358

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 848.6667
$ws.Range("J5").Value = 1404.6666
$ws.Range("L5").Value = 4213.9998
$ws.Range("N5").Value = -4437.9998

$ws.Range("H135").Value = 848.6667
$ws.Range("J135").Value = 1404.6666
$ws.Range("L135").Value = 12641.9994
$ws.Range("N135").Value = -17711.9994

$ws.Range("H136").Value = 3790
$ws.Range("I136").Value = 3386.6667
$ws.Range("K136").Value = 10160.0001
$ws.Range("M136").Value = -5060.000100000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1576.909
$ws.Range("I97").Value = 862.8
$ws.Range("J97").Value = 2172
$ws.Range("K97").Value = 862.8
$ws.Range("L97").Value = 2172
$ws.Range("M97").Value = -366.8
$ws.Range("N97").Value = -3164

$ws.Range("H99").Value = 7500
$ws.Range("I99").Value = 8250
$ws.Range("K99").Value = 8250
$ws.Range("M99").Value = -6004

$ws.Range("H113").Value = 2965.3635
$ws.Range("I113").Value = 2520
$ws.Range("J113").Value = 3499.8
$ws.Range("K113").Value = 2520
$ws.Range("L113").Value = 3499.8
$ws.Range("M113").Value = -350
$ws.Range("N113").Value = -7839.8

$ws.Range("H126").Value = 9800
$ws.Range("I126").Value = 3500
$ws.Range("K126").Value = 10500
$ws.Range("M126").Value = -8030

$ws.Range("H132").Value = 3231
$ws.Range("I132").Value = 4724.5
$ws.Range("K132").Value = 14173.5
$ws.Range("M132").Value = -11643.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3921.6
$ws.Range("I68").Value = 3861
$ws.Range("J68").Value = 4164
$ws.Range("K68").Value = 3861
$ws.Range("L68").Value = 4164
$ws.Range("M68").Value = -3112
$ws.Range("N68").Value = -5662

$ws.Range("H71").Value = 3921.6
$ws.Range("I71").Value = 3861
$ws.Range("J71").Value = 4164
$ws.Range("K71").Value = 19305
$ws.Range("L71").Value = 20820
$ws.Range("M71").Value = -15561
$ws.Range("N71").Value = -28308

$ws.Range("H93").Value = 12985.117
$ws.Range("I93").Value = 818.92
$ws.Range("J93").Value = 46780.11
$ws.Range("K93").Value = 818.92
$ws.Range("L93").Value = 46780.11
$ws.Range("M93").Value = 429.08
$ws.Range("N93").Value = -49276.11

$ws.Range("H100").Value = 253352.38
$ws.Range("I100").Value = 402230.47
$ws.Range("K100").Value = 402230.47
$ws.Range("M100").Value = -401689.47

$ws.Range("H136").Value = 76929280
$ws.Range("I136").Value = 4622
$ws.Range("J136").Value = 250009740
$ws.Range("K136").Value = 13866
$ws.Range("L136").Value = 750029220
$ws.Range("M136").Value = -11316
$ws.Range("N136").Value = -750034320

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 2007
$ws.Range("I113").Value = 662.5
$ws.Range("K113").Value = 1987.5
$ws.Range("M113").Value = 182.5

